# 22 Mayıs 2020 verileri eklendi
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# New data row for 2020-05-22
$row = 72
$ws.Cells.Item($row, 1).Value = 43973
$ws.Cells.Item($row, 2).Value = 37507
$ws.Cells.Item($row, 3).Value = 952
$ws.Cells.Item($row, 4).Value = 27
$ws.Cells.Item($row, 5).Value = 1121

# Expand the table (ListObject) to include the new row
$table = $ws.ListObjects.Item("Table3")
$table.Resize($ws.Range("A1:E72"))

# Update the active selection to match the recorded cursor position
$ws.Range("B70").Select()
